# Applies the "ajout des infos dans les contrats 25 aout 2025" edit:
# The bank's representative changes from
#   "Monsieur El Hadji Mamadou FAYE, son Directeur Général"
# to
#   "Madame Jenny MVOU, son Directeur Général Adjointe"
# (paragraph that starts "COFINA Gabon SA est représentée par ...").

$d = $word.ActiveDocument

# 1) Swap the bold proper name "El Hadji Mamadou FAYE" -> "Jenny MVOU".
#    Find/Replace on an existing run preserves that run's character
#    formatting (bold), matching the diff's bold "Jenny MVOU" run.
$rName = $d.Content
$rName.Find.Execute(
    "El Hadji Mamadou FAYE", $true, $false, $false, $false, $false,
    $true, 1, $false, "Jenny MVOU", 2)

# 2) Drop the "Monsieur " civility from the non-bold lead-in run.
$rMonsieur = $d.Content
$rMonsieur.Find.Execute(
    "est représentée par Monsieur ", $true, $false, $false, $false, $false,
    $true, 1, $false, "est représentée par ", 2)

# 3) Insert the new (non-bold) civility "Madame " right before the bold
#    name run; InsertBefore picks up the formatting of the preceding
#    (non-bold) text, so it lands as its own non-bold run exactly like
#    the target markup.
$rJenny = $d.Content
$rJenny.Find.Execute("Jenny MVOU", $true)
$rJenny.InsertBefore("Madame ")

# 4) Update the job title: "son Directeur Général, " -> "son Directeur
#    Général Adjointe, ".
$rTitle = $d.Content
$rTitle.Find.Execute(
    "son Directeur Général, ", $true, $false, $false, $false, $false,
    $true, 1, $false, "son Directeur Général Adjointe, ", 2)
